# ToDo.xlsx refactor:
#   - "Complete Renderer refactoring" (row 2) is finished -> remove the row entirely.
#   - "Deprecate Vector3..." (old row 4, now row 3) is replaced by the new task
#     "Move to pre-compiled shaders" with an estimate of 5.
#   - All rows below shift up by one; the two review comments (originally on
#     B13 and B16) need to stay attached to the same logical tasks, which are
#     now on B12 and B15.
#   - The sheet's last-used selection moves to B20.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the text of the two existing comments before touching rows, since
# they are anchored to cells that are about to move.
$shaderCommentText     = $ws.Range("B13").Comment.Text()
$degenerateCommentText = $ws.Range("B16").Comment.Text()

# Comments don't automatically follow a row delete, so drop them now and
# re-create them afterwards at their new locations.
$ws.Range("B13").Comment.Delete()
$ws.Range("B16").Comment.Delete()

# Row 2 ("Engine" / "Complete Renderer refactoring" / 10) is done -> delete it.
# Everything below (rows 3-19) shifts up to rows 2-18.
$ws.Rows.Item(2).Delete()

# The row that used to be row 4 ("Deprecate Vector3...", estimate 10) is now
# row 3; update it in place to reflect the new task.
$ws.Range("B3").Value = "Move to pre-compiled shaders"
$ws.Range("C3").Value = 5

# Re-attach the two comments to their shifted cells (old B13 -> B12,
# old B16 -> B15).
[void]$ws.Range("B12").AddComment($shaderCommentText)
[void]$ws.Range("B15").AddComment($degenerateCommentText)

# Match the final selection left behind in the saved workbook.
[void]$ws.Range("B20").Select()
